$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 30: Nmap scan-target entry, appended after the existing last row (29).
$ws.Range("A30").Value = 'Nmap'
$ws.Range("B30").Value = 'scan target'
$ws.Range("C30").Value = 'The format of nmap is: $nmap [Scan types] [Options] {target_ip_address}
$namp -sS 127.0.0.1
$nmap -sS 127.*.*.1-255
$nmap -sX -p 20,30,40,8080 127.0.0.1'

# Match the formatting of the row above it (font/wrap/alignment) rather than
# building a brand-new style.
$ws.Range("A29:C29").Copy()
$ws.Range("A30:C30").PasteSpecial(-4122)

# The wrapped 4-line description needs a taller row, same as other multi-line
# entries elsewhere in the sheet (e.g. row 7).
$ws.Rows.Item(30).RowHeight = 51.75

# Leave the new row selected, matching where the cursor ends up after typing
# the last entry and moving one cell to the right.
$null = $ws.Range("C31").Select()
